$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 0
$ws.Range("S2").Value = 68377
$ws.Range("T2").Value = 5373.441499999999
$ws.Range("U2").Value = 45000
$ws.Range("B3").Value = 12500
$ws.Range("I3").Value = 10000
$ws.Range("S3").Value = 65754
$ws.Range("T3").Value = 5250.357
$ws.Range("U3").Value = 22500
$ws.Range("S4").Value = 63452
$ws.Range("T4").Value = 5171.886999999999
$ws.Range("S5").Value = 61778
$ws.Range("T5").Value = 5107.980499999999
$ws.Range("S6").Value = 63701
$ws.Range("T6").Value = 5104.393
$ws.Range("S7").Value = 66921
$ws.Range("T7").Value = 5230.476999999999
$ws.Range("S8").Value = 65559
$ws.Range("T8").Value = 5902.879499999999
$ws.Range("P9").Value = 0
$ws.Range("S9").Value = 77302
$ws.Range("T9").Value = 7037.001999999999
$ws.Range("U9").Value = 45000
$ws.Range("S10").Value = 94353
$ws.Range("T10").Value = 8592.440500000001
$ws.Range("S11").Value = 84696
$ws.Range("T11").Value = 14095.0215
$ws.Range("S12").Value = 105351
$ws.Range("T12").Value = 15828.8445
$ws.Range("V12").Value = 5997.667411111111
$ws.Range("W12").Value = 5.478505529195085
$ws.Range("S13").Value = 106616
$ws.Range("T13").Value = 15449.287
$ws.Range("S14").Value = 105965
$ws.Range("T14").Value = 15803.2595
$ws.Range("T15").Value = 15862.028
$ws.Range("T16").Value = 16032.3205
$ws.Range("T17").Value = 16559.4275
$ws.Range("T18").Value = 16683.57249999999
$ws.Range("T19").Value = 16252.82049999999
$ws.Range("T20").Value = 15546.958
$ws.Range("T21").Value = 13748.4165
$ws.Range("T22").Value = 12030.7215
$ws.Range("T23").Value = 9319.645999999999
$ws.Range("T24").Value = 6308.333499999995
$ws.Range("P25").Value = 12000
$ws.Range("T25").Value = 5687.261999999999
$ws.Range("U25").Value = 57000

$ws.Range("S15").ClearContents()
$ws.Range("V15").ClearContents()
$ws.Range("W15").ClearContents()
